# Remove the row for ions that do not permeate ("98, 780, 1073"),
# shifting subsequent rows up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Delete()
